$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J62").Value = 3778
$ws.Range("L62").Value = 3778
$ws.Range("N62").Value = -5026
$ws.Range("H63").Value = 60000
$ws.Range("J63").Value = 60000
$ws.Range("L63").Value = 60000
$ws.Range("N63").Value = -61248
$ws.Range("J65").Value = 3778
$ws.Range("L65").Value = 18890
$ws.Range("N65").Value = -25130
$ws.Range("H66").Value = 60000
$ws.Range("J66").Value = 60000
$ws.Range("L66").Value = 180000
$ws.Range("N66").Value = -186240
$ws.Range("H98").Value = 856.1
$ws.Range("I98").Value = 864.8889
$ws.Range("K98").Value = 864.8889
$ws.Range("M98").Value = 633.1111
$ws.Range("H100").Value = 3150.75
$ws.Range("I100").Value = 1103.25
$ws.Range("K100").Value = 1103.25
$ws.Range("M100").Value = -562.25
$ws.Range("H116").Value = 17423.092
$ws.Range("J116").Value = 19801.5
$ws.Range("L116").Value = 19801.5
$ws.Range("N116").Value = -26685.5
$ws.Range("H122").Value = 856.1
$ws.Range("I122").Value = 864.8889
$ws.Range("K122").Value = 2594.6667
$ws.Range("M122").Value = -144.6667000000002
$ws.Range("H138").Value = 2672.7966
$ws.Range("J138").Value = 3106.0645
$ws.Range("L138").Value = 9318.193499999999
$ws.Range("N138").Value = -19598.1935

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1281.5778
$ws.Range("I32").Value = 913.4761999999999
$ws.Range("K32").Value = 913.4761999999999
$ws.Range("M32").Value = -626.4761999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 13223
$ws.Range("I62").Value = 4127.6
$ws.Range("K62").Value = 4127.6
$ws.Range("M62").Value = -3503.6
$ws.Range("H65").Value = 13223
$ws.Range("I65").Value = 4127.6
$ws.Range("K65").Value = 20638
$ws.Range("M65").Value = -17518
$ws.Range("H92").Value = 19130.2
$ws.Range("J92").Value = 19130.2
$ws.Range("L92").Value = 19130.2
$ws.Range("N92").Value = -24122.2
$ws.Range("H94").Value = 1975
$ws.Range("I94").Value = 1937.5
$ws.Range("J94").Value = 2005
$ws.Range("K94").Value = 1937.5
$ws.Range("L94").Value = 2005
$ws.Range("M94").Value = -1486.5
$ws.Range("N94").Value = -2907
$ws.Range("H103").Value = 38755.75
$ws.Range("I103").Value = 10012
$ws.Range("K103").Value = 10012
$ws.Range("M103").Value = -8840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1463432.4
$ws.Range("J11").Value = 261.2
$ws.Range("L11").Value = 783.5999999999999
$ws.Range("N11").Value = -1063.6
$ws.Range("H63").Value = 8895.25
$ws.Range("I63").Value = 8896.5
$ws.Range("J63").Value = 8894
$ws.Range("K63").Value = 26689.5
$ws.Range("L63").Value = 26682
$ws.Range("M63").Value = -25940.5
$ws.Range("N63").Value = -28180
$ws.Range("H66").Value = 8895.25
$ws.Range("I66").Value = 8896.5
$ws.Range("J66").Value = 8894
$ws.Range("K66").Value = 80068.5
$ws.Range("L66").Value = 80046
$ws.Range("M66").Value = -76324.5
$ws.Range("N66").Value = -87534
$ws.Range("H70").Value = 14998.75
$ws.Range("I70").Value = 9997.5
$ws.Range("K70").Value = 29992.5
$ws.Range("M70").Value = -29677.5
$ws.Range("H73").Value = 14998.75
$ws.Range("I73").Value = 9997.5
$ws.Range("K73").Value = 29992.5
$ws.Range("M73").Value = -28900.5
$ws.Range("H108").Value = 523.5
$ws.Range("I108").Value = 523.5
$ws.Range("K108").Value = 1570.5
$ws.Range("M108").Value = 1309.5
$ws.Range("H114").Value = 2316.625
$ws.Range("J114").Value = 3689.8
$ws.Range("L114").Value = 11069.4
$ws.Range("N114").Value = -17577.4
$ws.Range("H115").Value = 5996.3335
$ws.Range("I115").Value = 1994.5
$ws.Range("K115").Value = 5983.5
$ws.Range("M115").Value = -4808.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7362.561
$ws.Range("I70").Value = 7615.931
$ws.Range("K70").Value = 7615.931
$ws.Range("M70").Value = -7345.931
$ws.Range("H73").Value = 7362.561
$ws.Range("I73").Value = 7615.931
$ws.Range("K73").Value = 7615.931
$ws.Range("M73").Value = -6679.931
$ws.Range("H80").Value = 14447.375
$ws.Range("I80").Value = 1895.3334
$ws.Range("J80").Value = 21978.6
$ws.Range("K80").Value = 1895.3334
$ws.Range("L80").Value = 21978.6
$ws.Range("M80").Value = -897.3334
$ws.Range("N80").Value = -23974.6
$ws.Range("H83").Value = 14447.375
$ws.Range("I83").Value = 1895.3334
$ws.Range("J83").Value = 21978.6
$ws.Range("K83").Value = 9476.666999999999
$ws.Range("L83").Value = 109893
$ws.Range("M83").Value = -4484.666999999999
$ws.Range("N83").Value = -119877
$ws.Range("H92").Value = 10712.833
$ws.Range("I92").Value = 2728
$ws.Range("J92").Value = 11438.728
$ws.Range("K92").Value = 2728
$ws.Range("L92").Value = 11438.728
$ws.Range("M92").Value = -856
$ws.Range("N92").Value = -15182.728
$ws.Range("H132").Value = 1778.6
$ws.Range("I132").Value = 1481
$ws.Range("K132").Value = 4443
$ws.Range("M132").Value = -1913

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H22").Value = 1939.2307
$ws.Range("I22").Value = 868.3333
$ws.Range("J22").Value = 2857.1428
$ws.Range("K22").Value = 868.3333
$ws.Range("L22").Value = 2857.1428
$ws.Range("M22").Value = -573.3333
$ws.Range("N22").Value = -3447.1428
$ws.Range("H27").Value = 1939.2307
$ws.Range("I27").Value = 868.3333
$ws.Range("J27").Value = 2857.1428
$ws.Range("K27").Value = 868.3333
$ws.Range("L27").Value = 2857.1428
$ws.Range("M27").Value = -761.3333
$ws.Range("N27").Value = -3071.1428
$ws.Range("H55").Value = 431.2353
$ws.Range("I55").Value = 227.1
$ws.Range("K55").Value = 227.1
$ws.Range("M55").Value = -54.09999999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 87809.664
$ws.Range("J46").Value = 87809.664
$ws.Range("L46").Value = 87809.664
$ws.Range("N46").Value = -88271.664
$ws.Range("H100").Value = 1882.6666
$ws.Range("I100").Value = 431.33334
$ws.Range("K100").Value = 862.66668
$ws.Range("M100").Value = -321.66668
$ws.Range("H134").Value = 87809.664
$ws.Range("J134").Value = 87809.664
$ws.Range("L134").Value = 263428.992
$ws.Range("N134").Value = -268498.992
